# commit 21012025 F3 mappings
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 16 new rows above the current row 2, pushing the existing
# data rows (currently 2-7) down to rows 18-23.
$ws.Range("A2:A17").EntireRow.Insert()

# Populate the newly inserted rows with the new Object/Node mappings.
$ws.Range("A2").Value = "D00123500000000 AP5008"
$ws.Range("B2").Value = "P000001013"

$ws.Range("A3").Value = "D00118900000000 Ultra-High Purity Colloidal Silica PD Su"
$ws.Range("B3").Value = "P000001013"

$ws.Range("A4").Value = "D00123000000000 Optiplane 2260"
$ws.Range("B4").Value = "P000001013"

$ws.Range("A5").Value = "DP000010000000 Cu3886 low dishing low defect Cu Bulk Sl"
$ws.Range("B5").Value = "P000001016"

$ws.Range("A6").Value = "D00123300000000 Str Partnership IMEC TFM Share"
$ws.Range("B6").Value = "P000001012"

$ws.Range("A7").Value = "D00123200000000 Adv Pkg CMP (1501-50 BOOST)"
$ws.Range("B7").Value = "P000001013"

$ws.Range("A8").Value = "D00115200000000 Ultra-High Purity Colloidal Silica"
$ws.Range("B8").Value = "P000001013"

$ws.Range("A9").Value = "DP000005000000 DP1284 high-rate Cu bulk slurry for TI"
$ws.Range("B9").Value = "P000001016"

$ws.Range("A10").Value = "RD210008 Advanced particle characterization"
$ws.Range("B10").Value = "P000001016"

$ws.Range("A11").Value = "RD242007 HPD8700 for Micron HBM CMP"
$ws.Range("B11").Value = "P000001016"

$ws.Range("A12").Value = "CPHOENX_01010000604520 "
$ws.Range("B12").Value = "L000010145"

$ws.Range("A13").Value = "CPHOENX_01010000604540 "
$ws.Range("B13").Value = "L000010145"

$ws.Range("A14").Value = "CPHOENX_01010000604530 "
$ws.Range("B14").Value = "L000010145"

$ws.Range("A15").Value = "CORAERP_MMOR1975.4147 "
$ws.Range("B15").Value = "L000009822"

$ws.Range("A16").Value = "ITEMPEU_001300003829 "
$ws.Range("B16").Value = "L000013174"

$ws.Range("A17").Value = "CORAERP_MMOR1975.SMQA "
$ws.Range("B17").Value = "L000006414"
